$wb = $excel.ActiveWorkbook

# --- About sheet: add Notes section (rows 5-13) ---
$about = $wb.Worksheets.Item("About")

$about.Range("A5").Value = "Notes:"
$about.Range("A5").Font.Bold = $true

$about.Range("B5").Value  = "In the US, many coal plants are subject to rules requiring them "
$about.Range("B6").Value  = "to retrofit to meet enviromental guidelines. This requires"
$about.Range("B7").Value  = "a one time investment decision for plant owners. Because we don't"
$about.Range("B8").Value  = "track individual plants in the model, we calibrate the share of forward"
$about.Range("B9").Value  = "costs that must be recovered to represent the additional revenue that is needed to "
$about.Range("B10").Value = "save and pay for these one time investments and apply this across the distribution"
$about.Range("B11").Value = "of plant types. Calibration is done by comparing model results against other sources,"
$about.Range("B12").Value = "including Rhodium's ClimateDeck and EIA's Annual Energy Outlook and Electric "
$about.Range("B13").Value = "Power Monthly."

# --- SoFCtMbCtPR sheet: update calibrated value for "hard coal" ---
$data = $wb.Worksheets.Item("SoFCtMbCtPR ")

$data.Range("B2").Value = 2.5

$data.Range("B3").Select() | Out-Null

# Re-activate the About sheet last so it remains the selected tab
$about.Range("D31").Select() | Out-Null
